$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Restructure the "optimization_parameters" sheet (rows 8-17) ---
#
# Before:
#   8  Model                  | Sigmoid
#   9  estimate_params        | 1
#   10 make_graphs            | 0
#   11 fix_P                  | 1
#   12 fix_b                  | 1
#   13 expression_timepoints  | 0.4 | 0.8 | 1.2 | 1.6
#   14 Strain                 | wt  | dcin5
#   15 Sheet                  | 3   | 4
#   16 Deletion               | 0   | 3
#   17 simulation_timepoints  | 0 .. 2 (step 0.1)
#
# After:
#   8  production_function    | Sigmoid
#   9  L_curve                | 0
#   10 estimate_params        | 1
#   11 make_graphs            | 0
#   12 fix_P                  | 1
#   13 fix_b                  | 1
#   14 expression_timepoints  | 0.4 | 0.8 | 1.2 | 1.6
#   15 Strain                 | wt  | dcin5
#   16 Sheet                  | 3   | 4
#   17 simulation_timepoints  | 0 .. 2 (step 0.1)
#
# i.e. a new "L_curve" row is inserted right after the "Model"/"production_function"
# row, and the old "Deletion" row is removed - net row count is unchanged.

# 1) Insert a fresh row at position 9 - this pushes the old rows 9-17 down to 10-18.
$ws.Rows.Item(9).Insert()

# 2) Remove the old "Deletion" row, which now lives at row 17 (after the shift above).
$ws.Rows.Item(17).Delete()

# 3) Rename row 8's label from "Model" to "production_function" (value stays "Sigmoid").
$ws.Range("A8").Value = "production_function"

# 4) Populate the newly inserted row 9 with the "L_curve" parameter.
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value2 = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# --- View-state updates ---
# The active tab moves from "network_weights" (index 5) to "optimization_parameters"
# (index 6), and the selection on "optimization_parameters" moves to the full last row.
$ws.Activate()
$ws.Range("A17:XFD17").Select()
